$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 206, shifting existing rows 206:251 down to 207:252
$ws.Rows.Item(206).Insert()

# Populate the new row 206 with data (columns A-J copied from the surrounding dataset,
# columns D, K-T set to the new record's values)
$ws.Cells.Item(206, 1).Value = 1
$ws.Cells.Item(206, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(206, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(206, 4).Value = 44736
$ws.Cells.Item(206, 5).Value = 15
$ws.Cells.Item(206, 6).Value = "Fruta"
$ws.Cells.Item(206, 7).Value = 100102
$ws.Cells.Item(206, 8).Value = "Cítricos"
$ws.Cells.Item(206, 9).Value = 100102003
$ws.Cells.Item(206, 10).Value = "Limón"
$ws.Cells.Item(206, 11).Value = "Sin especificar"
$ws.Cells.Item(206, 12).Value = "3a amarillo"
$ws.Cells.Item(206, 13).Value = 300
$ws.Cells.Item(206, 14).Value = 8500
$ws.Cells.Item(206, 15).Value = 9000
$ws.Cells.Item(206, 16).Value = 8750
$ws.Cells.Item(206, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(206, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(206, 19).Value = 438
$ws.Cells.Item(206, 20).Value = 20
